$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the header style/format from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for I2:J78 (I0 / IF columns), row-aligned with existing data rows 2-78
$data = @(
    @(6, 7),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(5, 5)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
